# Rename task worksheets: replace spaces with underscores (and drop the
# extra internal space, e.g. "Relaxing Music1" -> "RelaxingMusic1")
$wb = $excel.ActiveWorkbook

$renameMap = @{
    "Task1 Relaxing Music1"  = "Task1_RelaxingMusic1"
    "Task2 Relaxing Pic1"    = "Task2_RelaxingPic1"
    "Task3 EMA1"             = "Task3_EMA1"
    "Task4 IAPS"             = "Task4_IAPS"
    "Task5 StressPic1"       = "Task5_StressPic1"
    "Task6 EMA2"             = "Task6_EMA2"
    "Task7 Relaxing Music2"  = "Task7_RelaxingMusic2"
    "Task8 Relaxing Pic2"    = "Task8_RelaxingPic2"
    "Task9 EMA3"             = "Task9_EMA3"
    "Task10 Biking"          = "Task10_Biking"
    "Task11 Relaxing Music3" = "Task11_RelaxingMusic3"
    "Task12 NeutralPic"      = "Task12_NeutralPic"
    "Task13 EMA4"            = "Task13_EMA4"
    "Task14 MentalMath"      = "Task14_MentalMath"
    "Task15 Stroop"          = "Task15_Stroop"
    "Task16 StressPic2"      = "Task16_StressPic2"
    "Task17 EMA5"            = "Task17_EMA5"
    "Task18 March"           = "Task18_March"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renameMap.ContainsKey($oldName)) {
        $ws.Name = $renameMap[$oldName]
    }
}
